$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.432888150215149
$ws.Range("B1").Value = 2.59496021270752
$ws.Range("C1").Value = 4.001776218414307
$ws.Range("D1").Value = 4.059993267059326
$ws.Range("E1").Value = 2.386115550994873
